$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(43, 8).Value = 4417.1  # H43
$ws.Cells.Item(43, 10).Value = 4417.1  # J43
$ws.Cells.Item(43, 12).Value = 4417.1  # L43
$ws.Cells.Item(43, 14).Value = -4555.1  # N43
$ws.Cells.Item(58, 8).Value = 1554.0333  # H58
$ws.Cells.Item(58, 9).Value = 407.7143  # I58
$ws.Cells.Item(58, 10).Value = 1902.9131  # J58
$ws.Cells.Item(58, 11).Value = 1223.1429  # K58
$ws.Cells.Item(58, 12).Value = 5708.7393  # L58
$ws.Cells.Item(58, 13).Value = -1073.1429  # M58
$ws.Cells.Item(58, 14).Value = -6008.7393  # N58
$ws.Cells.Item(74, 8).Value = 4400.6665  # H74
$ws.Cells.Item(74, 9).Value = 4400.6665  # I74
$ws.Cells.Item(74, 11).Value = 4400.6665  # K74
$ws.Cells.Item(74, 13).Value = -3464.6665  # M74
$ws.Cells.Item(77, 8).Value = 4400.6665  # H77
$ws.Cells.Item(77, 9).Value = 4400.6665  # I77
$ws.Cells.Item(77, 11).Value = 22003.3325  # K77
$ws.Cells.Item(77, 13).Value = -17323.3325  # M77
$ws.Cells.Item(100, 8).Value = 2265.9333  # H100
$ws.Cells.Item(100, 9).Value = 1876.2307  # I100
$ws.Cells.Item(100, 10).Value = 4799  # J100
$ws.Cells.Item(100, 11).Value = 1876.2307  # K100
$ws.Cells.Item(100, 12).Value = 4799  # L100
$ws.Cells.Item(100, 13).Value = -1335.2307  # M100
$ws.Cells.Item(100, 14).Value = -5881  # N100
$ws.Cells.Item(107, 8).Value = 728.4074000000001  # H107
$ws.Cells.Item(107, 9).Value = 728.4074000000001  # I107
$ws.Cells.Item(107, 11).Value = 728.4074000000001  # K107
$ws.Cells.Item(107, 13).Value = 1191.5926  # M107
$ws.Cells.Item(116, 8).Value = 336172.44  # H116
$ws.Cells.Item(116, 9).Value = 601940.8  # I116
$ws.Cells.Item(116, 10).Value = 3962  # J116
$ws.Cells.Item(116, 11).Value = 601940.8  # K116
$ws.Cells.Item(116, 12).Value = 3962  # L116
$ws.Cells.Item(116, 13).Value = -598498.8  # M116
$ws.Cells.Item(116, 14).Value = -10846  # N116
$ws.Cells.Item(131, 8).Value = 3261.7  # H131
$ws.Cells.Item(131, 9).Value = 645.3333  # I131
$ws.Cells.Item(131, 10).Value = 7186.25  # J131
$ws.Cells.Item(131, 11).Value = 1935.9999  # K131
$ws.Cells.Item(131, 12).Value = 21558.75  # L131
$ws.Cells.Item(131, 13).Value = 3104.0001  # M131
$ws.Cells.Item(131, 14).Value = -31638.75  # N131
$ws.Cells.Item(135, 8).Value = 28734.658  # H135
$ws.Cells.Item(135, 9).Value = 968.61536  # I135
$ws.Cells.Item(135, 10).Value = 88894.414  # J135
$ws.Cells.Item(135, 11).Value = 8717.53824  # K135
$ws.Cells.Item(135, 12).Value = 800049.726  # L135
$ws.Cells.Item(135, 13).Value = -6182.53824  # M135
$ws.Cells.Item(135, 14).Value = -805119.726  # N135
$ws.Cells.Item(138, 8).Value = 5160.2666  # H138
$ws.Cells.Item(138, 10).Value = 5102.516  # J138
$ws.Cells.Item(138, 12).Value = 15307.548  # L138
$ws.Cells.Item(138, 14).Value = -25587.548  # N138

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 14012.618  # H32
$ws.Cells.Item(32, 9).Value = 8279.146000000001  # I32
$ws.Cells.Item(32, 10).Value = 30803.5  # J32
$ws.Cells.Item(32, 11).Value = 8279.146000000001  # K32
$ws.Cells.Item(32, 12).Value = 30803.5  # L32
$ws.Cells.Item(32, 13).Value = -7992.146000000001  # M32
$ws.Cells.Item(32, 14).Value = -31377.5  # N32
$ws.Cells.Item(45, 8).Value = 1116699.4  # H45
$ws.Cells.Item(45, 10).Value = 6286.625  # J45
$ws.Cells.Item(45, 12).Value = 6286.625  # L45
$ws.Cells.Item(45, 14).Value = -7040.625  # N45
$ws.Cells.Item(97, 8).Value = 2019.95  # H97
$ws.Cells.Item(97, 9).Value = 966.73334  # I97
$ws.Cells.Item(97, 11).Value = 966.73334  # K97
$ws.Cells.Item(97, 13).Value = -470.73334  # M97
$ws.Cells.Item(110, 8).Value = 1680.2354  # H110
$ws.Cells.Item(110, 9).Value = 1639.9354  # I110
$ws.Cells.Item(110, 11).Value = 1639.9354  # K110
$ws.Cells.Item(110, 13).Value = 405.0645999999999  # M110
$ws.Cells.Item(132, 8).Value = 21349.309  # H132
$ws.Cells.Item(132, 9).Value = 21571.844  # I132
$ws.Cells.Item(132, 11).Value = 64715.53200000001  # K132
$ws.Cells.Item(132, 13).Value = -62185.53200000001  # M132

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(99, 8).Value = 93148.45  # H99
$ws.Cells.Item(99, 10).Value = 3066.3333  # J99
$ws.Cells.Item(99, 12).Value = 3066.3333  # L99
$ws.Cells.Item(99, 14).Value = -6062.3333  # N99

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(107, 8).Value = 250462.42  # H107
$ws.Cells.Item(107, 9).Value = 250462.42  # I107
$ws.Cells.Item(107, 11).Value = 250462.42  # K107
$ws.Cells.Item(107, 13).Value = -248542.42  # M107
$ws.Cells.Item(132, 8).Value = 3503.3076  # H132
$ws.Cells.Item(132, 9).Value = 3503.3076  # I132
$ws.Cells.Item(132, 11).Value = 10509.9228  # K132
$ws.Cells.Item(132, 13).Value = -7979.9228  # M132
$ws.Cells.Item(134, 8).Value = 25898.363  # H134
$ws.Cells.Item(134, 9).Value = 28143.35  # I134
$ws.Cells.Item(134, 10).Value = 3448.5  # J134
$ws.Cells.Item(134, 11).Value = 84430.04999999999  # K134
$ws.Cells.Item(134, 12).Value = 10345.5  # L134
$ws.Cells.Item(134, 13).Value = -81895.04999999999  # M134
$ws.Cells.Item(134, 14).Value = -15415.5  # N134

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(37, 8).Value = 61963.11  # H37
$ws.Cells.Item(37, 10).Value = 61963.11  # J37
$ws.Cells.Item(37, 12).Value = 185889.33  # L37
$ws.Cells.Item(37, 14).Value = -186113.33  # N37
$ws.Cells.Item(80, 8).Value = 4293.727  # H80
$ws.Cells.Item(80, 9).Value = 4489  # I80
$ws.Cells.Item(80, 10).Value = 4274.2  # J80
$ws.Cells.Item(80, 11).Value = 13467  # K80
$ws.Cells.Item(80, 12).Value = 12822.6  # L80
$ws.Cells.Item(80, 13).Value = -12531  # M80
$ws.Cells.Item(80, 14).Value = -14694.6  # N80
$ws.Cells.Item(83, 8).Value = 4293.727  # H83
$ws.Cells.Item(83, 9).Value = 4489  # I83
$ws.Cells.Item(83, 10).Value = 4274.2  # J83
$ws.Cells.Item(83, 11).Value = 40401  # K83
$ws.Cells.Item(83, 12).Value = 38467.8  # L83
$ws.Cells.Item(83, 13).Value = -35721  # M83
$ws.Cells.Item(83, 14).Value = -47827.8  # N83
$ws.Cells.Item(92, 8).Value = 952.7619  # H92
$ws.Cells.Item(92, 9).Value = 584  # I92
$ws.Cells.Item(92, 10).Value = 1444.4445  # J92
$ws.Cells.Item(92, 11).Value = 1752  # K92
$ws.Cells.Item(92, 12).Value = 4333.333500000001  # L92
$ws.Cells.Item(92, 13).Value = -504  # M92
$ws.Cells.Item(92, 14).Value = -6829.333500000001  # N92
$ws.Cells.Item(122, 8).Value = 1085.6923  # H122
$ws.Cells.Item(122, 9).Value = 911.5  # I122
$ws.Cells.Item(122, 11).Value = 8203.5  # K122
$ws.Cells.Item(122, 13).Value = -5753.5  # M122

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(2, 8).Value = 152.61539  # H2
$ws.Cells.Item(2, 10).Value = 268  # J2
$ws.Cells.Item(2, 12).Value = 268  # L2
$ws.Cells.Item(2, 14).Value = -494  # N2
$ws.Cells.Item(132, 8).Value = 34128.883  # H132
$ws.Cells.Item(132, 9).Value = 37302.71  # I132
$ws.Cells.Item(132, 11).Value = 111908.13  # K132
$ws.Cells.Item(132, 13).Value = -109378.13  # M132
$ws.Cells.Item(139, 8).Value = 69339.664  # H139
$ws.Cells.Item(139, 10).Value = 69339.664  # J139
$ws.Cells.Item(139, 12).Value = 69339.664  # L139
$ws.Cells.Item(139, 14).Value = -79619.664  # N139

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 806.913  # H55
$ws.Cells.Item(55, 9).Value = 552.4545000000001  # I55
$ws.Cells.Item(55, 11).Value = 552.4545000000001  # K55
$ws.Cells.Item(55, 13).Value = -379.4545000000001  # M55
$ws.Cells.Item(93, 8).Value = 1685.1052  # H93
$ws.Cells.Item(93, 9).Value = 1534.5333  # I93
$ws.Cells.Item(93, 11).Value = 1534.5333  # K93
$ws.Cells.Item(93, 13).Value = -286.5333000000001  # M93
$ws.Cells.Item(100, 8).Value = 4598.2  # H100
$ws.Cells.Item(100, 9).Value = 5246  # I100
$ws.Cells.Item(100, 10).Value = 4166.3335  # J100
$ws.Cells.Item(100, 11).Value = 5246  # K100
$ws.Cells.Item(100, 12).Value = 4166.3335  # L100
$ws.Cells.Item(100, 13).Value = -4705  # M100
$ws.Cells.Item(100, 14).Value = -5248.3335  # N100
$ws.Cells.Item(132, 8).Value = 51446.92  # H132
$ws.Cells.Item(132, 9).Value = 63097.15  # I132
$ws.Cells.Item(132, 10).Value = 4846  # J132
$ws.Cells.Item(132, 11).Value = 189291.45  # K132
$ws.Cells.Item(132, 12).Value = 14538  # L132
$ws.Cells.Item(132, 13).Value = -186761.45  # M132
$ws.Cells.Item(132, 14).Value = -19598  # N132
$ws.Cells.Item(136, 8).Value = 4450.1333  # H136
$ws.Cells.Item(136, 9).Value = 3329.2  # I136
$ws.Cells.Item(136, 11).Value = 9987.599999999999  # K136
$ws.Cells.Item(136, 13).Value = -7437.599999999999  # M136

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(96, 8).Value = 5778.8887  # H96
$ws.Cells.Item(96, 9).Value = 3802.2  # I96
$ws.Cells.Item(96, 11).Value = 3802.2  # K96
$ws.Cells.Item(96, 13).Value = -2429.2  # M96
$ws.Cells.Item(107, 8).Value = 691.9  # H107
$ws.Cells.Item(107, 9).Value = 445.7143  # I107
$ws.Cells.Item(107, 10).Value = 1266.3334  # J107
$ws.Cells.Item(107, 11).Value = 1337.1429  # K107
$ws.Cells.Item(107, 12).Value = 3799.0002  # L107
$ws.Cells.Item(107, 13).Value = 582.8571000000002  # M107
$ws.Cells.Item(107, 14).Value = -7639.0002  # N107
$ws.Cells.Item(132, 8).Value = 35641.76  # H132
$ws.Cells.Item(132, 9).Value = 39562.734  # I132
$ws.Cells.Item(132, 11).Value = 118688.202  # K132
$ws.Cells.Item(132, 13).Value = -116158.202  # M132
$ws.Cells.Item(133, 8).Value = 84999  # H133
$ws.Cells.Item(133, 10).Value = 84999  # J133
$ws.Cells.Item(133, 12).Value = 84999  # L133
$ws.Cells.Item(133, 14).Value = -95119  # N133
$ws.Cells.Item(136, 8).Value = 5719391.5  # H136
$ws.Cells.Item(136, 9).Value = 597532.3  # I136
$ws.Cells.Item(136, 11).Value = 1792596.9  # K136
$ws.Cells.Item(136, 13).Value = -1790046.9  # M136
